$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking values in columns D and E so Excel
# does not auto-convert them to numbers (the source file stores every data
# cell as text/inline string).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '64.255.73'
$ws.Range('E2').Value = '  +5.39%  '
$ws.Range('D3').Value = '2.744.63'
$ws.Range('E3').Value = '  +2.84%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '582.75'
$ws.Range('E5').Value = '  +2.77%  '
$ws.Range('D6').Value = '158.54'
$ws.Range('E6').Value = '  +9.70%  '
$ws.Range('E7').Value = '  -0.25%  '
$ws.Range('D8').Value = '0.612'
$ws.Range('E8').Value = '  +2.02%  '
$ws.Range('D9').Value = '2.770.29'
$ws.Range('E9').Value = '  +3.82%  '
$ws.Range('D10').Value = '6.77'
$ws.Range('E10').Value = '  +2.60%  '
$ws.Range('E11').Value = '  +6.41%  '
$ws.Range('E12').Value = '  +3.84%  '
$ws.Range('E13').Value = '  +2.08%  '
$ws.Range('D14').Value = '3.234.99'
$ws.Range('E14').Value = '  +3.08%  '
$ws.Range('D15').Value = '27.33'
$ws.Range('E15').Value = '  +4.50%  '
$ws.Range('D16').Value = '64.130.43'
$ws.Range('E16').Value = '  +5.20%  '
$ws.Range('E17').Value = '  +7.73%  '
$ws.Range('D18').Value = '2.765.07'
$ws.Range('E18').Value = '  +3.85%  '
$ws.Range('E19').Value = '  +4.85%  '
$ws.Range('D20').Value = '4.97'
$ws.Range('E20').Value = '  +4.83%  '
$ws.Range('D21').Value = '364.58'
$ws.Range('E21').Value = '  +3.83%  '
$ws.Range('E22').Value = '  +2.78%  '
$ws.Range('D23').Value = '0.998'
$ws.Range('E23').Value = '  -0.22%  '
$ws.Range('E24').Value = '  +1.31%  '
$ws.Range('D25').Value = '67.25'
$ws.Range('E25').Value = '  +5.08%  '
$ws.Range('E26').Value = '  +5.87%  '
$ws.Range('E27').Value = '  +5.40%  '
$ws.Range('E28').Value = '  -0.27%  '
$ws.Range('D29').Value = '0.0₃0918'
$ws.Range('E29').Value = '  +13.23%  '
$ws.Range('E30').Value = '  +1.27%  '
$ws.Range('E31').Value = '  +5.77%  '
$ws.Range('E32').Value = '  +19.53%  '
$ws.Range('D33').Value = '174.89'
$ws.Range('E33').Value = '  +6.86%  '
$ws.Range('D34').Value = '20.77'
$ws.Range('E34').Value = '  +4.24%  '
$ws.Range('E35').Value = '  -0.21%  '
$ws.Range('D36').Value = '4.93'
$ws.Range('E36').Value = '  +7.03%  '
$ws.Range('E37').Value = '  +9.50%  '
$ws.Range('E38').Value = '  +10.34%  '
$ws.Range('D39').Value = '1.02'
$ws.Range('E39').Value = '  +11.13%  '
$ws.Range('E40').Value = '  +5.54%  '
$ws.Range('D41').Value = '342.88'
$ws.Range('E41').Value = '  +1.01%  '
$ws.Range('D42').Value = '39.47'
$ws.Range('E42').Value = '  +2.41%  '
$ws.Range('D43').Value = '5.89'
$ws.Range('E43').Value = '  +13.03%  '
$ws.Range('D44').Value = '22.20'
$ws.Range('E44').Value = '  +8.98%  '
$ws.Range('D45').Value = '22.23'
$ws.Range('E45').Value = '  +7.81%  '
$ws.Range('D46').Value = '0.0603'
$ws.Range('E46').Value = '  +6.77%  '
$ws.Range('E47').Value = '  +4.77%  '
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').Value = '0.0262'
$ws.Range('E48').Value = '  +4.90%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').Value = '139.05'
$ws.Range('E49').Value = '  +4.45%  '
$ws.Range('E50').Value = '  +2.37%  '

# Restore the default (General) formatting/style so the cells keep the same
# style index as the original workbook.
$ws.Range("D2:E51").ClearFormats()

